$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 corresponds to the "Experimental" property; set its value (column B) to the
# literal text "false" (leading apostrophe forces text, not a boolean/number).
$ws.Cells.Item(7, 2).Value = "'false"
# Re-apply the same formatting used by the rest of the column so the cell keeps
# its original "wrap text / bordered" look instead of the quote-prefix style.
$ws.Cells.Item(6, 2).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4122)

# Row 8 corresponds to the "Date" property; refresh its timestamp
$ws.Cells.Item(8, 2).Value = "2025-10-03T16:37:46+01:00"
